$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1 headers ----
# Existing header cells (A1:I1 already carry the bold/centered/bordered style) -- just update the text
$ws.Range("B1").Value = "Alexis Rainey"
$ws.Range("C1").Value = "Balduzzi"
$ws.Range("D1").Value = "Burns"
$ws.Range("E1").Value = "Curley"
$ws.Range("F1").Value = "Doyle"
$ws.Range("G1").Value = "Espona"
$ws.Range("H1").Value = "Ferriolo"
$ws.Range("I1").Value = "Hackman"

# New header cells (J1:W1) -- set the text, then copy the format from A1 (which already
# has the bold/centered/bordered header style) so the new cells reuse the same style
# index instead of Excel fabricating near-duplicate style entries.
$ws.Range("J1").Value = "Holzman"
$ws.Range("A1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("K1").Value = "Hughes"
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("L1").Value = "Johnson"
$ws.Range("A1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$ws.Range("M1").Value = "McCann"
$ws.Range("A1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("N1").Value = "McFadden"
$ws.Range("A1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("O1").Value = "Medico"
$ws.Range("A1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$ws.Range("P1").Value = "Myers"
$ws.Range("A1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("Q1").Value = "Pla"
$ws.Range("A1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)
$ws.Range("R1").Value = "Reilly"
$ws.Range("A1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$ws.Range("S1").Value = "Rodrigo"
$ws.Range("A1").Copy()
$ws.Range("S1").PasteSpecial(-4122)
$ws.Range("T1").Value = "Streib"
$ws.Range("A1").Copy()
$ws.Range("T1").PasteSpecial(-4122)
$ws.Range("U1").Value = "Tollaksen"
$ws.Range("A1").Copy()
$ws.Range("U1").PasteSpecial(-4122)
$ws.Range("V1").Value = "Wasyliw"
$ws.Range("A1").Copy()
$ws.Range("V1").PasteSpecial(-4122)
$ws.Range("W1").Value = "Yanovich"
$ws.Range("A1").Copy()
$ws.Range("W1").PasteSpecial(-4122)

# ---- Data rows 2-6 ----
# A leading apostrophe forces Excel to store an empty *text* cell (matching the blank
# inlineStr cells in the target) rather than simply deleting the cell. ClearFormats()
# afterwards strips the "quote prefix" flag Excel would otherwise remember for that cell,
# so the cell ends up with no style attribute at all, matching the target exactly.
# Row 2
$ws.Range("B2").Value = 15.8203
$ws.Range("C2").Value = 14.0731
$ws.Range("D2").Value = 17.2049
$ws.Range("E2").Value = "'"
$ws.Range("E2").ClearFormats()
$ws.Range("F2").Value = 20.2307
$ws.Range("G2").Value = 17.0505
$ws.Range("H2").Value = 11.9108
$ws.Range("I2").Value = 14.4191
$ws.Range("J2").Value = 13.5592
$ws.Range("K2").Value = "'"
$ws.Range("K2").ClearFormats()
$ws.Range("L2").Value = 13.2302
$ws.Range("M2").Value = 16.1707
$ws.Range("N2").Value = 18.7414
$ws.Range("O2").Value = 16.5156
$ws.Range("P2").Value = 15.2876
$ws.Range("Q2").Value = 11.9207
$ws.Range("R2").Value = "'"
$ws.Range("R2").ClearFormats()
$ws.Range("S2").Value = 13.9856
$ws.Range("T2").Value = "'"
$ws.Range("T2").ClearFormats()
$ws.Range("U2").Value = 18.7996
$ws.Range("V2").Value = 15.8072
$ws.Range("W2").Value = 14.3708

# Row 3
$ws.Range("B3").Value = 9.7939
$ws.Range("C3").Value = 7.1131
$ws.Range("D3").Value = 6.2012
$ws.Range("E3").Value = 3.8481
$ws.Range("F3").Value = 8.0408
$ws.Range("G3").Value = 6.297
$ws.Range("H3").Value = 5.4767
$ws.Range("I3").Value = 7.5679
$ws.Range("J3").Value = 6.5123
$ws.Range("K3").Value = 4.2128
$ws.Range("L3").Value = 9.7423
$ws.Range("M3").Value = 8.3328
$ws.Range("N3").Value = 9.2733
$ws.Range("O3").Value = 11.082
$ws.Range("P3").Value = 6.7916
$ws.Range("Q3").Value = 4.405
$ws.Range("R3").Value = 6.0587
$ws.Range("S3").Value = 8.2017
$ws.Range("T3").Value = "'"
$ws.Range("T3").ClearFormats()
$ws.Range("U3").Value = 8.3067
$ws.Range("V3").Value = 8.5824
$ws.Range("W3").Value = 8.6075

# Row 4
$ws.Range("B4").Value = 10.8107
$ws.Range("C4").Value = 8.6637
$ws.Range("D4").Value = 10.3986
$ws.Range("E4").Value = 3.9678
$ws.Range("F4").Value = 15.3591
$ws.Range("G4").Value = 9.376
$ws.Range("H4").Value = 14.2764
$ws.Range("I4").Value = 12.3825
$ws.Range("J4").Value = 14.7108
$ws.Range("K4").Value = 8.2395
$ws.Range("L4").Value = 15.986
$ws.Range("M4").Value = 14.5628
$ws.Range("N4").Value = 13.2121
$ws.Range("O4").Value = 18.4323
$ws.Range("P4").Value = 14.2141
$ws.Range("Q4").Value = 9.5451
$ws.Range("R4").Value = 11.1606
$ws.Range("S4").Value = 10.9264
$ws.Range("T4").Value = "'"
$ws.Range("T4").ClearFormats()
$ws.Range("U4").Value = 15.499
$ws.Range("V4").Value = 12.1438
$ws.Range("W4").Value = "'"
$ws.Range("W4").ClearFormats()

# Row 5
$ws.Range("B5").Value = 11.1859
$ws.Range("C5").Value = 8.2695
$ws.Range("D5").Value = "'"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = 4.0078
$ws.Range("F5").Value = 11.6853
$ws.Range("G5").Value = 10.977
$ws.Range("H5").Value = 7.8411
$ws.Range("I5").Value = 8.2457
$ws.Range("J5").Value = 8.219
$ws.Range("K5").Value = 5.4288
$ws.Range("L5").Value = 10.3625
$ws.Range("M5").Value = 9.3619
$ws.Range("N5").Value = 13.0407
$ws.Range("O5").Value = 11.9271
$ws.Range("P5").Value = "'"
$ws.Range("P5").ClearFormats()
$ws.Range("Q5").Value = 9.6564
$ws.Range("R5").Value = 8.1462
$ws.Range("S5").Value = 9.1576
$ws.Range("T5").Value = 3.5187
$ws.Range("U5").Value = 11.8114
$ws.Range("V5").Value = 9.8252
$ws.Range("W5").Value = "'"
$ws.Range("W5").ClearFormats()

# Row 6
$ws.Range("B6").Value = 14.8602
$ws.Range("C6").Value = 10.404
$ws.Range("D6").Value = "'"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = 3.8997
$ws.Range("F6").Value = 16.0503
$ws.Range("G6").Value = 13.3825
$ws.Range("H6").Value = 12.1056
$ws.Range("I6").Value = 15.3125
$ws.Range("J6").Value = 11.7627
$ws.Range("K6").Value = 4.6309
$ws.Range("L6").Value = 11.7589
$ws.Range("M6").Value = 8.9039
$ws.Range("N6").Value = 13.7584
$ws.Range("O6").Value = 19.2116
$ws.Range("P6").Value = "'"
$ws.Range("P6").ClearFormats()
$ws.Range("Q6").Value = 8.4769
$ws.Range("R6").Value = 9.9943
$ws.Range("S6").Value = 9.7444
$ws.Range("T6").Value = 3.6805
$ws.Range("U6").Value = "'"
$ws.Range("U6").ClearFormats()
$ws.Range("V6").Value = 10.5224
$ws.Range("W6").Value = "'"
$ws.Range("W6").ClearFormats()

